$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.440.67'
$ws.Range('E2').Value = '  -0.26%  '
$ws.Range('D3').Value = '1.566.46'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('E5').Value = '  +0.11%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '287.84'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.92%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3729'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.80%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '48.21'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -3.60%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3317'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.30%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.131'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.63%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07466'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.33%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.002'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.21%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '20.63'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -3.03%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.934'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.68%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.902'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.40%  '
$ws.Range('D16').Value = '1.561.83'
$ws.Range('E16').Value = '  -0.80%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001112'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.09%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.06756'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.57%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '87.70'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -3.14%  '
$ws.Range('E20').Value = '  +0.16%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.337'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.52%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '16.41'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '12.07'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.01%  '
$ws.Range('D24').Value = '22.432.66'
$ws.Range('E24').Value = '  -0.24%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.383'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.77%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.563'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -3.44%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '152.87'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.03%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '19.66'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.95%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.009'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.93%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '123.99'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.79%  '
$ws.Range('D31').Value = '1.738.55'
$ws.Range('E31').Value = '  -0.68%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.051'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.82%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.011'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.22%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.116'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.74%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '9.625'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.53%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.08282'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.45%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02455'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.21%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2274'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.25%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.06380'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -2.44%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.350'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -1.72%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.286'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -4.60%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '11.24'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.77%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.6264'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.23%  '
$ws.Range('E44').Value = '  +0.14%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.77'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.53%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.6103'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +3.70%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.767'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.88%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.039'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.55%  '
$ws.Range('E49').Value = '  -1.86%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.209'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.37%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.07221'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.26%  '
